$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the per-row custom "zone" labels in column I (rows 2..42) with a
# repeating 5-zone pattern: Зона 1, Зона 2, Зона 3, Зона 4, Зона 5, Зона 1, ...
$zoneNames = @("Зона 1", "Зона 2", "Зона 3", "Зона 4", "Зона 5")

for ($row = 2; $row -le 42; $row++) {
    $zoneIndex = ($row - 2) % 5
    $ws.Cells.Item($row, 9).Value = $zoneNames[$zoneIndex]
}
